$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Rename the inline picture "image2.png" -> "image1.png" for the Pearson
# Edexcel logo that appears in the footers, and "image1.jpg" -> "image2.jpg"
# for the BTEC logo in the header. InlineShape has no writable Name
# property (same as real Word), so each picture is momentarily converted
# to a floating Shape (which does expose Name), renamed, then converted
# back to an inline shape so the drawing stays wp:inline exactly as before.

function Rename-InlinePicture($range, $newName) {
    $shapeCount = $range.InlineShapes.Count
    if ($shapeCount -ge 1) {
        $inlineShape = $range.InlineShapes.Item(1)
        $shape = $inlineShape.ConvertToShape()
        $shape.Name = $newName
        $shape.ConvertToInlineShape() | Out-Null
    }
}

for ($i = 1; $i -le $sec.Headers.Count; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists) {
        foreach ($ishp in $hdr.Range.InlineShapes) {
            if ($ishp.AlternativeText -eq "BTec_Logo-Orange") {
                Rename-InlinePicture $hdr.Range "image2.jpg"
            }
        }
    }
}

for ($i = 1; $i -le $sec.Footers.Count; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists) {
        foreach ($ishp in $ftr.Range.InlineShapes) {
            if ($ishp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                Rename-InlinePicture $ftr.Range "image1.png"
            }
        }
    }
}

Write-Output "Renamed inline pictures."
